$d = $word.ActiveDocument

# Locate the paragraph that ends the first quote block:
# "Das bringt Abschied und Ankunft und neues Berühren.”1"
$target = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text -like "*Das bringt Abschied und Ankunft und neues Ber*hren*") {
        $target = $para
        break
    }
}

# Insert a blank paragraph right after the target paragraph.
$target.Range.InsertParagraphAfter()

# Insert a second blank paragraph after that one; this will become the
# new "vielleicht ein kleiner Bezug ..." note paragraph.
$blank = $d.Paragraphs($target.Index + 1)
$blank.Range.InsertParagraphAfter()

# The paragraph that will hold the arrow marker.
$notePara = $d.Paragraphs($target.Index + 2)
$noteRange = $notePara.Range
$noteRange.InsertAfter([char]0x2192 + " ")

# Split the text into two distinct runs by temporarily putting the rest
# of the sentence into its own paragraph, then merging the paragraph
# break away -- this keeps "→ " and the following sentence as two
# separate <w:r> runs instead of Word silently coalescing them into one.
$noteRange.InsertParagraphAfter()
$secondPara = $d.Paragraphs($notePara.Index + 1)
$secondPara.Range.InsertAfter("vielleicht ein kleiner Bezug zur Biographie von Rainer Schell mit dem man das Zitat in den Text einbauen kann.")

$mergePos = $notePara.Range.End - 1
$mergeRange = $d.Range($mergePos, $mergePos + 1)
$mergeRange.Delete()
